$d = $word.ActiveDocument

# 1. "Welcome to Polymorphism Visualisation" -> add trailing period
$d.Content.Find.Execute("Welcome to Polymorphism Visualisation", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Welcome to Polymorphism Visualisation.", 2)

# 2. "...assigned to the variable." -> "...assigned to that variable."
#    (only the intro paragraph occurrence; use unique longer context to avoid
#    touching the similar phrase later in the document)
$d.Content.Find.Execute("the type of the object assigned to the variable.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the type of the object assigned to that variable.", 2)

# 3. "different shapes represents a different class" -> "different shapes represent a different class"
$d.Content.Find.Execute("different shapes represents a different class", $true, $false, $false, $false, $false,
                         $true, 1, $false, "different shapes represent a different class", 2)

# 4. "in to the variable at run time" -> "into the variable at run time"
$d.Content.Find.Execute("passed in to the variable", $true, $false, $false, $false, $false,
                         $true, 1, $false, "passed into the variable", 2)

# 5. Merge "5. " + "Clear button" runs into a single run with the same text.
$d.Content.Find.Execute("5. Clear button", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5. Clear button", 2)

# 6. Merge "8" + ". " runs into a single run with the same text.
$d.Content.Find.Execute("8. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8. ", 2)
